# Insert two new rows (DADA2_Species, DADA2_Taxonomy) after the CustomNBC row (row 4)
# and before the Kraken2_0.0 row (currently row 5). This shifts all subsequent
# rows down by 2, matching the new dimension A1:F15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows starting at row 5, pushing existing rows 5.. down to 7..
$ws.Range("A5:F6").EntireRow.Insert()

# Fill in the new DADA2_Species row (row 5)
$ws.Range("A5").Value = "DADA2_Species"
$ws.Range("B5").Value = 0.02
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

# Fill in the new DADA2_Taxonomy row (row 6)
$ws.Range("A6").Value = "DADA2_Taxonomy"
$ws.Range("B6").Value = 0.17
$ws.Range("C6").Value = 0.28
$ws.Range("D6").Value = 0.25
$ws.Range("E6").Value = 0.31
$ws.Range("F6").Value = 0.22
